$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.794.11'
$ws.Range("E2").Value = '  +4.82%  '
$ws.Range("D3").Value = '1.610.25'
$ws.Range("E3").Value = '  +3.65%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").Formula = "=""213.75"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").Formula = "=""0.516"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +6.94%  '
$ws.Range("D7").Formula = "=""0.995"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("D8").Formula = "=""26.97"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +11.82%  '
$ws.Range("E9").Value = '  +3.19%  '
$ws.Range("E10").Value = '  +2.48%  '
$ws.Range("D11").Formula = "=""0.0912"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").Value = '1.840.24'
$ws.Range("E12").Value = '  +3.62%  '
$ws.Range("D13").Value = '1.606.72'
$ws.Range("D14").Value = '29.825.42'
$ws.Range("E14").Value = '  +4.91%  '
$ws.Range("D15").Formula = "=""0.539"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +5.73%  '
$ws.Range("D16").Formula = "=""3.76"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  +3.77%  '
$ws.Range("D17").Formula = "=""244.69"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +7.02%  '
$ws.Range("E18").Value = '  +3.75%  '
$ws.Range("E19").Value = '  +3.75%  '
$ws.Range("D20").Value = '0.0₃0693'
$ws.Range("E20").Value = '  +3.10%  '
$ws.Range("D21").Formula = "=""0.996"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").Formula = "=""4.05"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +4.12%  '
$ws.Range("D23").Formula = "=""9.25"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +3.81%  '
$ws.Range("D24").Formula = "=""2.10"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +4.22%  '
$ws.Range("D25").Formula = "=""155.52"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +2.82%  '
$ws.Range("D26").Formula = "=""15.34"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +4.21%  '
$ws.Range("E27").Value = '  +5.50%  '
$ws.Range("E28").Value = '  +2.66%  '
$ws.Range("D29").Formula = "=""0.996"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("D30").Formula = "=""0.0473"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Formula = "=""3.25"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("D33").Value = '1.446.66'
$ws.Range("E33").Value = '  +4.53%  '
$ws.Range("E34").Value = '  +3.54%  '
$ws.Range("D35").Formula = "=""1.04"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("D36").Formula = "=""2.83"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +9.95%  '
$ws.Range("E37").Value = '  +2.47%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  +3.24%  '
$ws.Range("E40").Value = '  +4.88%  '
$ws.Range("D41").Formula = "=""55.60"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +30.15%  '
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").Formula = "=""0.795"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +2.95%  '
$ws.Range("D45").Formula = "=""0.0468"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").Formula = "=""66.54"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +7.76%  '
$ws.Range("D47").Formula = "=""5.31"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '1.751.74'
$ws.Range("E48").Value = '  +3.87%  '
$ws.Range("D49").Formula = "=""86.96"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +2.26%  '
$ws.Range("D50").Formula = "=""0.836"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -4.20%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Formula = "=""0.0522"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +2.24%  '
$excel.CutCopyMode = 0
